$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "2021-10-05 13:40:45.680935",
    "2021-10-05 13:40:45.680947",
    "2021-10-05 13:40:45.680951",
    "2021-10-05 13:40:45.680954",
    "2021-10-05 13:40:45.680957",
    "2021-10-05 13:40:45.680960",
    "2021-10-05 13:40:45.680963",
    "2021-10-05 13:40:45.680966",
    "2021-10-05 13:40:45.680970",
    "2021-10-05 13:40:45.680973",
    "2021-10-05 13:40:45.680976",
    "2021-10-05 13:40:45.680979",
    "2021-10-05 13:40:45.680982",
    "2021-10-05 13:40:45.680985",
    "2021-10-05 13:40:45.680988",
    "2021-10-05 13:40:45.680990",
    "2021-10-05 13:40:45.680994",
    "2021-10-05 13:40:45.680997",
    "2021-10-05 13:40:45.681000",
    "2021-10-05 13:40:45.681003",
    "2021-10-05 13:40:45.681006",
    "2021-10-05 13:40:45.681009",
    "2021-10-05 13:40:45.681012",
    "2021-10-05 13:40:45.681014",
    "2021-10-05 13:40:45.681018",
    "2021-10-05 13:40:45.681021"
)

# New header cell F1: give it the same formatting as the existing header
# cells (bold font + border, centered) by copying E1's format onto F1,
# then write the header text.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Fill in the per-row "time_taken" metadata column for every data row.
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
